{"js": "// Update the division-problem worksheet table: each cell's equation text is\n// replaced with a new equation. Cells are addressed by (row, column)\n// position (not just by matching old text) because some source equations\n// repeat (e.g. \"47\u00f77=\" appears in two different cells but maps to two\n// different replacements), so a document-wide text search/replace would be\n// ambiguous. Replacing via a cell-scoped search().insertText(\"Replace\")\n// only rewrites the <w:t> run content and keeps the existing run/paragraph\n// formatting (font, size, alignment) intact.\nconst replacements = [\n  { row: 0, col: 0, oldText: \"71\u00f77=\", newText: \"11\u00f78=\" },\n  { row: 0, col: 1, oldText: \"44\u00f72=\", newText: \"17\u00f77=\" },\n  { row: 0, col: 2, oldText: \"25\u00f76=\", newText: \"17\u00f74=\" },\n  { row: 0, col: 3, oldText: \"47\u00f77=\", newText: \"81\u00f75=\" },\n  { row: 0, col: 4, oldText: \"58\u00f77=\", newText: \"90\u00f77=\" },\n\n  { row: 4, col: 0, oldText: \"32\u00f74=\", newText: \"97\u00f75=\" },\n  { row: 4, col: 1, oldText: \"41\u00f75=\", newText: \"31\u00f74=\" },\n  { row: 4, col: 2, oldText: \"47\u00f77=\", newText: \"94\u00f76=\" },\n  { row: 4, col: 3, oldText: \"43\u00f76=\", newText: \"27\u00f77=\" },\n  { row: 4, col: 4, oldText: \"59\u00f78=\", newText: \"78\u00f72=\" },\n\n  { row: 8, col: 0, oldText: \"14\u00f75=\", newText: \"36\u00f76=\" },\n  { row: 8, col: 1, oldText: \"40\u00f74=\", newText: \"71\u00f75=\" },\n  { row: 8, col: 2, oldText: \"39\u00f77=\", newText: \"32\u00f78=\" },\n  { row: 8, col: 3, oldText: \"12\u00f75=\", newText: \"48\u00f76=\" },\n  { row: 8, col: 4, oldText: \"19\u00f73=\", newText: \"26\u00f77=\" },\n\n  { row: 12, col: 0, oldText: \"70\u00f77=\", newText: \"43\u00f79=\" },\n  { row: 12, col: 1, oldText: \"24\u00f73=\", newText: \"89\u00f74=\" },\n  { row: 12, col: 2, oldText: \"28\u00f73=\", newText: \"50\u00f75=\" },\n  { row: 12, col: 3, oldText: \"81\u00f72=\", newText: \"71\u00f78=\" },\n  { row: 12, col: 4, oldText: \"67\u00f76=\", newText: \"89\u00f75=\" },\n\n  { row: 16, col: 0, oldText: \"42\u00f78=\", newText: \"43\u00f76=\" },\n  { row: 16, col: 1, oldText: \"59\u00f77=\", newText: \"73\u00f75=\" },\n  { row: 16, col: 2, oldText: \"94\u00f79=\", newText: \"63\u00f76=\" },\n  { row: 16, col: 3, oldText: \"36\u00f79=\", newText: \"68\u00f79=\" },\n  { row: 16, col: 4, oldText: \"88\u00f78=\", newText: \"48\u00f73=\" },\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nfor (const { row, col, oldText, newText } of replacements) {\n  const cell = table.getCell(row, col);\n  const hits = cell.body.search(oldText, { matchCase: true });\n  hits.load(\"items\");\n  await context.sync();\n\n  if (hits.items.length > 0) {\n    hits.items[0].insertText(newText, \"Replace\");\n  } else {\n    // Fallback (shouldn't happen for this document): overwrite the whole\n    // cell text if the expected equation wasn't found verbatim.\n    cell.body.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the division-problem worksheet table: each cell's equation text is\n# replaced with a new equation. Cells are addressed by (row, column)\n# position (1-based, matching Word's Table.Cell(row, col) indexing) because\n# some source equations repeat (e.g. \"47\u00f77=\" appears in two different\n# cells but maps to two different replacements), so a document-wide\n# Find/Replace would be ambiguous (and this host's Find.Execute searches the\n# whole story regardless of the Range it is scoped to). Assigning directly\n# to Cell.Range.Text only rewrites the run's text and leaves the existing\n# run/paragraph formatting (font, size, alignment) untouched.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$replacements = @(\n  @{ Row = 1; Col = 1; Old = \"71\u00f77=\"; New = \"11\u00f78=\" },\n  @{ Row = 1; Col = 2; Old = \"44\u00f72=\"; New = \"17\u00f77=\" },\n  @{ Row = 1; Col = 3; Old = \"25\u00f76=\"; New = \"17\u00f74=\" },\n  @{ Row = 1; Col = 4; Old = \"47\u00f77=\"; New = \"81\u00f75=\" },\n  @{ Row = 1; Col = 5; Old = \"58\u00f77=\"; New = \"90\u00f77=\" },\n\n  @{ Row = 5; Col = 1; Old = \"32\u00f74=\"; New = \"97\u00f75=\" },\n  @{ Row = 5; Col = 2; Old = \"41\u00f75=\"; New = \"31\u00f74=\" },\n  @{ Row = 5; Col = 3; Old = \"47\u00f77=\"; New = \"94\u00f76=\" },\n  @{ Row = 5; Col = 4; Old = \"43\u00f76=\"; New = \"27\u00f77=\" },\n  @{ Row = 5; Col = 5; Old = \"59\u00f78=\"; New = \"78\u00f72=\" },\n\n  @{ Row = 9; Col = 1; Old = \"14\u00f75=\"; New = \"36\u00f76=\" },\n  @{ Row = 9; Col = 2; Old = \"40\u00f74=\"; New = \"71\u00f75=\" },\n  @{ Row = 9; Col = 3; Old = \"39\u00f77=\"; New = \"32\u00f78=\" },\n  @{ Row = 9; Col = 4; Old = \"12\u00f75=\"; New = \"48\u00f76=\" },\n  @{ Row = 9; Col = 5; Old = \"19\u00f73=\"; New = \"26\u00f77=\" },\n\n  @{ Row = 13; Col = 1; Old = \"70\u00f77=\"; New = \"43\u00f79=\" },\n  @{ Row = 13; Col = 2; Old = \"24\u00f73=\"; New = \"89\u00f74=\" },\n  @{ Row = 13; Col = 3; Old = \"28\u00f73=\"; New = \"50\u00f75=\" },\n  @{ Row = 13; Col = 4; Old = \"81\u00f72=\"; New = \"71\u00f78=\" },\n  @{ Row = 13; Col = 5; Old = \"67\u00f76=\"; New = \"89\u00f75=\" },\n\n  @{ Row = 17; Col = 1; Old = \"42\u00f78=\"; New = \"43\u00f76=\" },\n  @{ Row = 17; Col = 2; Old = \"59\u00f77=\"; New = \"73\u00f75=\" },\n  @{ Row = 17; Col = 3; Old = \"94\u00f79=\"; New = \"63\u00f76=\" },\n  @{ Row = 17; Col = 4; Old = \"36\u00f79=\"; New = \"68\u00f79=\" },\n  @{ Row = 17; Col = 5; Old = \"88\u00f78=\"; New = \"48\u00f73=\" }\n)\n\nforeach ($item in $replacements) {\n  $cell = $t.Cell($item.Row, $item.Col)\n  $cell.Range.Text = $item.New\n}\n"}
